$wb = $excel.ActiveWorkbook

# --- "Stella 6.1" sheet: rename to "Stella", update row 2 values, move selection ---
$wsStella = $wb.Worksheets.Item(1)
$wsStella.Name = "Stella"
$wsStella.Range("B2").Value = "Red"
$wsStella.Range("C2").Value = "TEST"
$wsStella.Range("A2").Select()

# --- "TRAIL Neo 3" sheet: add row 2 data, make it the active/selected tab ---
$wsTrail = $wb.Worksheets.Item(2)
$wsTrail.Range("A2").Value = "L"
$wsTrail.Range("B2").Value = "White"
$wsTrail.Range("C2").Value = "test"
$wsTrail.Activate()
$wsTrail.Range("C2").Select()
